$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("ForecastValidations")
$wsCurrent  = $wb.Worksheets.Item("CurrentWeatherValidations")

# --- Content edits -------------------------------------------------------

# ForecastValidations: retype the city names to include a space.
$wsForecast.Range("B3").Value = "Aundh Camp"
$wsForecast.Range("B4").Value = "New Delhi"

# CurrentWeatherValidations: correct the latitude precision for Phoenix.
$wsCurrent.Range("E5").Value = 33.44838

# --- Selection / active-sheet state --------------------------------------
# End state: CurrentWeatherValidations keeps B3 selected (no longer the
# active tab), ForecastValidations becomes the active tab with B4 selected.

$wsCurrent.Range("B3").Select()

$wsForecast.Select()
$wsForecast.Range("B4").Select()
